# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures on the
# crypto tracking sheet, as produced by the scheduled GitHub Actions update.
#
# Note: a handful of "Price" values are plain decimals (e.g. 244.46) that
# Excel would otherwise auto-convert into numbers on assignment. Those are
# written with a leading apostrophe ('') to force them to stay as literal
# text, matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.692.03'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.960.61'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''244.46'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '''0.618'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = '''58.94'
$ws.Range("E7").Value = '  +1.81%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").Value = '''0.0806'
$ws.Range("E10").Value = '  -3.07%  '
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '''22.17'
$ws.Range("E12").Value = '  +3.25%  '
$ws.Range("D13").Value = '2.248.80'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '''0.823'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '''13.71'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").Value = '''5.27'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = '1.962.45'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '36.567.49'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '''69.66'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '0.0₃0860'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").Value = '''228.93'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '''5.08'
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("D25").Value = '''2.34'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '''9.31'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +13.54%  '
$ws.Range("D28").Value = '''160.41'
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("D29").Value = '''19.36'
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").Value = '''4.26'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = '''6.08'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("E37").Value = '  +4.36%  '
$ws.Range("D38").Value = '''3.39'
$ws.Range("E38").Value = '  +12.87%  '
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  +3.67%  '
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").Value = '1.358.61'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").Value = '''87.49'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("D50").Value = '2.139.48'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").Value = '''43.46'
$ws.Range("E51").Value = '  -4.81%  '
